$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(58, 8).Value = 193.5  # H58, old=149.71428
$ws.Cells.Item(62, 8).Value = 4064.3225  # H62, old=4064.3547
$ws.Cells.Item(62, 9).Value = 2769.16  # I62, old=2701.1155
$ws.Cells.Item(62, 10).Value = 9460.833000000001  # J62, old=11153.2
$ws.Cells.Item(62, 11).Value = 2769.16  # K62, old=2701.1155
$ws.Cells.Item(62, 12).Value = 9460.833000000001  # L62, old=11153.2
$ws.Cells.Item(62, 13).Value = -2145.16  # M62, old=-2077.1155
$ws.Cells.Item(62, 14).Value = -10708.833  # N62, old=-12401.2
$ws.Cells.Item(65, 8).Value = 4064.3225  # H65, old=4064.3547
$ws.Cells.Item(65, 9).Value = 2769.16  # I65, old=2701.1155
$ws.Cells.Item(65, 10).Value = 9460.833000000001  # J65, old=11153.2
$ws.Cells.Item(65, 11).Value = 13845.8  # K65, old=13505.5775
$ws.Cells.Item(65, 12).Value = 47304.165  # L65, old=55766
$ws.Cells.Item(65, 13).Value = -10725.8  # M65, old=-10385.5775
$ws.Cells.Item(65, 14).Value = -53544.165  # N65, old=-62006
$ws.Cells.Item(92, 8).Value = 1325.9375  # H92, old=1366.7742
$ws.Cells.Item(92, 9).Value = 1363.7742  # I92, old=1407.2333
$ws.Cells.Item(92, 11).Value = 1363.7742  # K92, old=1407.2333
$ws.Cells.Item(92, 13).Value = -115.7742000000001  # M92, old=-159.2333000000001
$ws.Cells.Item(97, 8).Value = 3026  # H97, old=3026.25
$ws.Cells.Item(97, 10).Value = 3368  # J97, old=3368.3333
$ws.Cells.Item(97, 12).Value = 10104  # L97, old=10104.9999
$ws.Cells.Item(97, 14).Value = -11096  # N97, old=-11096.9999
$ws.Cells.Item(132, 8).Value = 20001576  # H132, old=20409762
$ws.Cells.Item(132, 9).Value = 21278144  # I132, old=21740702
$ws.Cells.Item(132, 11).Value = 63834432  # K132, old=65222106
$ws.Cells.Item(132, 13).Value = -63831902  # M132, old=-65219576
$ws.Cells.Item(138, 8).Value = 3261.1729  # H138, old=3302.0854
$ws.Cells.Item(138, 9).Value = 2118.1035  # I138, old=2186.7407
$ws.Cells.Item(138, 10).Value = 3898.6538  # J138, old=3849.6182
$ws.Cells.Item(138, 11).Value = 6354.310500000001  # K138, old=6560.222099999999
$ws.Cells.Item(138, 12).Value = 11695.9614  # L138, old=11548.8546
$ws.Cells.Item(138, 13).Value = -1214.310500000001  # M138, old=-1420.222099999999
$ws.Cells.Item(138, 14).Value = -21975.9614  # N138, old=-21828.8546

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2222884  # H2, old=3472996.8
$ws.Cells.Item(2, 9).Value = 2222884  # I2, old=3472996.8
$ws.Cells.Item(2, 11).Value = 2222884  # K2, old=3472996.8
$ws.Cells.Item(2, 13).Value = -2222771  # M2, old=-3472883.8
$ws.Cells.Item(37, 8).Value = 39191.25  # H37, old=39222.5
$ws.Cells.Item(37, 10).Value = 44882.5  # J37, old=44945
$ws.Cells.Item(37, 12).Value = 44882.5  # L37, old=44945
$ws.Cells.Item(37, 14).Value = -45428.5  # N37, old=-45491
$ws.Cells.Item(61, 8).Value = 3903.1277  # H61, old=4080.8445
$ws.Cells.Item(61, 9).Value = 4089.2927  # I61, old=4142.927
$ws.Cells.Item(61, 10).Value = 2631  # J61, old=3444.5
$ws.Cells.Item(61, 11).Value = 4089.2927  # K61, old=4142.927
$ws.Cells.Item(61, 12).Value = 2631  # L61, old=3444.5
$ws.Cells.Item(61, 13).Value = -3877.2927  # M61, old=-3930.927
$ws.Cells.Item(61, 14).Value = -3055  # N61, old=-3868.5
$ws.Cells.Item(102, 8).Value = 5212849.5  # H102, old=4906179.5
$ws.Cells.Item(102, 9).Value = 9261495  # I102, old=7577764.5
$ws.Cells.Item(102, 10).Value = 7449.143  # J102, old=8274
$ws.Cells.Item(102, 11).Value = 9261495  # K102, old=7577764.5
$ws.Cells.Item(102, 12).Value = 7449.143  # L102, old=8274
$ws.Cells.Item(102, 13).Value = -9259873  # M102, old=-7576142.5
$ws.Cells.Item(102, 14).Value = -10693.143  # N102, old=-11518
$ws.Cells.Item(116, 8).Value = 2222884  # H116, old=3472996.8
$ws.Cells.Item(116, 9).Value = 2222884  # I116, old=3472996.8
$ws.Cells.Item(116, 11).Value = 2222884  # K116, old=3472996.8
$ws.Cells.Item(116, 13).Value = -2220590  # M116, old=-3470702.8
$ws.Cells.Item(136, 8).Value = 3903.1277  # H136, old=4080.8445
$ws.Cells.Item(136, 9).Value = 4089.2927  # I136, old=4142.927
$ws.Cells.Item(136, 10).Value = 2631  # J136, old=3444.5
$ws.Cells.Item(136, 11).Value = 12267.8781  # K136, old=12428.781
$ws.Cells.Item(136, 12).Value = 7893  # L136, old=10333.5
$ws.Cells.Item(136, 13).Value = -9717.8781  # M136, old=-9878.780999999999
$ws.Cells.Item(136, 14).Value = -12993  # N136, old=-15433.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2222884  # H3, old=3472996.8
$ws.Cells.Item(3, 9).Value = 2222884  # I3, old=3472996.8
$ws.Cells.Item(3, 11).Value = 2222884  # K3, old=3472996.8
$ws.Cells.Item(3, 13).Value = -2222770  # M3, old=-3472882.8
$ws.Cells.Item(22, 8).Value = 537.5  # H22, old=650
$ws.Cells.Item(22, 10).Value = 1000  # J22, old=916.6667
$ws.Cells.Item(22, 12).Value = 1000  # L22, old=916.6667
$ws.Cells.Item(22, 14).Value = -1346  # N22, old=-1262.6667
$ws.Cells.Item(54, 8).Value = 3500  # H54, old=6599.6665
$ws.Cells.Item(54, 9).Value = 650  # I54, old=400
$ws.Cells.Item(54, 10).Value = 9200  # J54, old=9699.5
$ws.Cells.Item(54, 11).Value = 650  # K54, old=400
$ws.Cells.Item(54, 12).Value = 9200  # L54, old=9699.5
$ws.Cells.Item(54, 13).Value = -166  # M54, old=84
$ws.Cells.Item(54, 14).Value = -10168  # N54, old=-10667.5
$ws.Cells.Item(80, 8).Value = 342.04544  # H80, old=330.3913
$ws.Cells.Item(80, 10).Value = 341.2  # J80, old=324.5
$ws.Cells.Item(80, 12).Value = 341.2  # L80, old=324.5
$ws.Cells.Item(80, 14).Value = -2337.2  # N80, old=-2320.5
$ws.Cells.Item(83, 8).Value = 342.04544  # H83, old=330.3913
$ws.Cells.Item(83, 10).Value = 341.2  # J83, old=324.5
$ws.Cells.Item(83, 12).Value = 1706  # L83, old=1622.5
$ws.Cells.Item(83, 14).Value = -11690  # N83, old=-11606.5
$ws.Cells.Item(86, 8).Value = 5569802.5  # H86, old=4177825
$ws.Cells.Item(86, 9).Value = 6676296.5  # I86, old=5007790
$ws.Cells.Item(86, 10).Value = 37331.668  # J86, old=28000
$ws.Cells.Item(86, 11).Value = 6676296.5  # K86, old=5007790
$ws.Cells.Item(86, 12).Value = 37331.668  # L86, old=28000
$ws.Cells.Item(86, 13).Value = -6675173.5  # M86, old=-5006667
$ws.Cells.Item(86, 14).Value = -39577.668  # N86, old=-30246
$ws.Cells.Item(89, 8).Value = 5569802.5  # H89, old=4177825
$ws.Cells.Item(89, 9).Value = 6676296.5  # I89, old=5007790
$ws.Cells.Item(89, 10).Value = 37331.668  # J89, old=28000
$ws.Cells.Item(89, 11).Value = 33381482.5  # K89, old=25038950
$ws.Cells.Item(89, 12).Value = 186658.34  # L89, old=140000
$ws.Cells.Item(89, 13).Value = -33375866.5  # M89, old=-25033334
$ws.Cells.Item(89, 14).Value = -197890.34  # N89, old=-151232
$ws.Cells.Item(94, 8).Value = 3854263.8  # H94, old=4550476
$ws.Cells.Item(94, 9).Value = 5264359  # I94, old=6668073
$ws.Cells.Item(94, 10).Value = 26862.572  # J94, old=12769
$ws.Cells.Item(94, 11).Value = 5264359  # K94, old=6668073
$ws.Cells.Item(94, 12).Value = 26862.572  # L94, old=12769
$ws.Cells.Item(94, 13).Value = -5263908  # M94, old=-6667622
$ws.Cells.Item(94, 14).Value = -27764.572  # N94, old=-13671
$ws.Cells.Item(107, 8).Value = 1986137.2  # H107, old=1833454.9
$ws.Cells.Item(107, 9).Value = 2552782.2  # I107, old=2464817.2
$ws.Cells.Item(107, 10).Value = 2879.625  # J107, old=2503.6
$ws.Cells.Item(107, 11).Value = 2552782.2  # K107, old=2464817.2
$ws.Cells.Item(107, 12).Value = 2879.625  # L107, old=2503.6
$ws.Cells.Item(107, 13).Value = -2550862.2  # M107, old=-2462897.2
$ws.Cells.Item(107, 14).Value = -6719.625  # N107, old=-6343.6
$ws.Cells.Item(134, 8).Value = 6503.021  # H134, old=6758.804
$ws.Cells.Item(134, 9).Value = 5761.878  # I134, old=6025.564
$ws.Cells.Item(134, 11).Value = 17285.634  # K134, old=18076.692
$ws.Cells.Item(134, 13).Value = -14750.634  # M134, old=-15541.692

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 5339.909  # H99, old=5429.5
$ws.Cells.Item(99, 9).Value = 5361  # I99, old=5666.6665
$ws.Cells.Item(99, 11).Value = 5361  # K99, old=5666.6665
$ws.Cells.Item(99, 13).Value = -3863  # M99, old=-4168.6665
$ws.Cells.Item(126, 8).Value = 5339.909  # H126, old=5429.5
$ws.Cells.Item(126, 9).Value = 5361  # I126, old=5666.6665
$ws.Cells.Item(126, 11).Value = 16083  # K126, old=16999.9995
$ws.Cells.Item(126, 13).Value = -13613  # M126, old=-14529.9995
$ws.Cells.Item(132, 8).Value = 1807.38  # H132, old=1840.7142
$ws.Cells.Item(132, 9).Value = 1714.9318  # I132, old=1737.9535
$ws.Cells.Item(132, 10).Value = 2485.3333  # J132, old=2577.1667
$ws.Cells.Item(132, 11).Value = 5144.7954  # K132, old=5213.860500000001
$ws.Cells.Item(132, 12).Value = 7455.999899999999  # L132, old=7731.500100000001
$ws.Cells.Item(132, 13).Value = -2614.7954  # M132, old=-2683.860500000001
$ws.Cells.Item(132, 14).Value = -12515.9999  # N132, old=-12791.5001
$ws.Cells.Item(134, 8).Value = 22469.072  # H134, old=23662.717
$ws.Cells.Item(134, 9).Value = 26458.098  # I134, old=27765.025
$ws.Cells.Item(134, 10).Value = 11565.733  # J134, old=12234.857
$ws.Cells.Item(134, 11).Value = 79374.29400000001  # K134, old=83295.07500000001
$ws.Cells.Item(134, 12).Value = 34697.199  # L134, old=36704.571
$ws.Cells.Item(134, 13).Value = -76839.29400000001  # M134, old=-80760.07500000001
$ws.Cells.Item(134, 14).Value = -39767.199  # N134, old=-41774.571

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(56, 8).Value = 16672960  # H56, old=15631304
$ws.Cells.Item(56, 9).Value = 16672960  # I56, old=15631304
$ws.Cells.Item(56, 11).Value = 16672960  # K56, old=15631304
$ws.Cells.Item(56, 13).Value = -16672430  # M56, old=-15630774
$ws.Cells.Item(93, 8).Value = 4999  # H93, old=6000
$ws.Cells.Item(93, 9).Value = 0  # I93, old=6000
$ws.Cells.Item(93, 10).Value = 4999  # J93, old=0
$ws.Cells.Item(93, 11).Value = 0  # K93, old=18000
$ws.Cells.Item(93, 12).Value = 14997  # L93, old=0
$ws.Cells.Item(93, 13).ClearContents()  # M93, old=-16128
$ws.Cells.Item(93, 14).Value = -18741  # N93, old=None
$ws.Cells.Item(94, 8).Value = 7456  # H94, old=4133.25
$ws.Cells.Item(94, 9).Value = 707.3333  # I94, old=634.4
$ws.Cells.Item(94, 10).Value = 9986.75  # J94, old=9964.666999999999
$ws.Cells.Item(94, 11).Value = 2121.9999  # K94, old=1903.2
$ws.Cells.Item(94, 12).Value = 29960.25  # L94, old=29894.001
$ws.Cells.Item(94, 13).Value = -1445.9999  # M94, old=-1227.2
$ws.Cells.Item(94, 14).Value = -31312.25  # N94, old=-31246.001
$ws.Cells.Item(114, 8).Value = 674.2105  # H114, old=526973.75
$ws.Cells.Item(114, 9).Value = 146.25  # I114, old=148.33333
$ws.Cells.Item(114, 10).Value = 815  # J114, old=625753.5
$ws.Cells.Item(114, 11).Value = 438.75  # K114, old=444.99999
$ws.Cells.Item(114, 12).Value = 2445  # L114, old=1877260.5
$ws.Cells.Item(114, 13).Value = 2815.25  # M114, old=2809.00001
$ws.Cells.Item(114, 14).Value = -8953  # N114, old=-1883768.5
$ws.Cells.Item(129, 8).Value = 963.5625  # H129, old=1043.25
$ws.Cells.Item(129, 9).Value = 793.6667  # I129, old=812.4167
$ws.Cells.Item(129, 10).Value = 1473.25  # J129, old=1735.75
$ws.Cells.Item(129, 11).Value = 2381.0001  # K129, old=2437.2501
$ws.Cells.Item(129, 12).Value = 4419.75  # L129, old=5207.25
$ws.Cells.Item(129, 13).Value = 2618.9999  # M129, old=2562.7499
$ws.Cells.Item(129, 14).Value = -14419.75  # N129, old=-15207.25

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 821894.5  # H97, old=768811.5
$ws.Cells.Item(97, 9).Value = 992956.3  # I97, old=916506.4399999999
$ws.Cells.Item(97, 11).Value = 992956.3  # K97, old=916506.4399999999
$ws.Cells.Item(97, 13).Value = -992460.3  # M97, old=-916010.4399999999
$ws.Cells.Item(126, 8).Value = 8052459.5  # H126, old=12881826
$ws.Cells.Item(126, 9).Value = 3790920.5  # I126, old=6496420.5
$ws.Cells.Item(126, 10).Value = 20837078  # J126, old=27781108
$ws.Cells.Item(126, 11).Value = 11372761.5  # K126, old=19489261.5
$ws.Cells.Item(126, 12).Value = 62511234  # L126, old=83343324
$ws.Cells.Item(126, 13).Value = -11370291.5  # M126, old=-19486791.5
$ws.Cells.Item(126, 14).Value = -62516174  # N126, old=-83348264

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 6197.4287  # H7, old=6229.2
$ws.Cells.Item(7, 9).Value = 4778.727  # I7, old=5386.625
$ws.Cells.Item(7, 10).Value = 11399.333  # J7, old=9599.5
$ws.Cells.Item(7, 11).Value = 4778.727  # K7, old=5386.625
$ws.Cells.Item(7, 12).Value = 11399.333  # L7, old=9599.5
$ws.Cells.Item(7, 13).Value = -4666.727  # M7, old=-5274.625
$ws.Cells.Item(7, 14).Value = -11623.333  # N7, old=-9823.5
$ws.Cells.Item(40, 8).Value = 8199.8125  # H40, old=7885.0586
$ws.Cells.Item(40, 9).Value = 7120  # I40, old=6731.727
$ws.Cells.Item(40, 11).Value = 7120  # K40, old=6731.727
$ws.Cells.Item(40, 13).Value = -6984  # M40, old=-6595.727
$ws.Cells.Item(55, 8).Value = 1092.3125  # H55, old=1159.8
$ws.Cells.Item(55, 9).Value = 1143.35  # I55, old=1261.5
$ws.Cells.Item(55, 11).Value = 1143.35  # K55, old=1261.5
$ws.Cells.Item(55, 13).Value = -970.3499999999999  # M55, old=-1088.5
$ws.Cells.Item(80, 8).Value = 54999.8  # H80, old=55499.75
$ws.Cells.Item(80, 10).Value = 56875  # J80, old=58166.668
$ws.Cells.Item(80, 12).Value = 56875  # L80, old=58166.668
$ws.Cells.Item(80, 14).Value = -59121  # N80, old=-60412.668
$ws.Cells.Item(83, 8).Value = 54999.8  # H83, old=55499.75
$ws.Cells.Item(83, 10).Value = 56875  # J83, old=58166.668
$ws.Cells.Item(83, 12).Value = 170625  # L83, old=174500.004
$ws.Cells.Item(83, 14).Value = -181857  # N83, old=-185732.004
$ws.Cells.Item(100, 8).Value = 36296.805  # H100, old=37406.734
$ws.Cells.Item(100, 10).Value = 128562.125  # J100, old=146499.72
$ws.Cells.Item(100, 12).Value = 128562.125  # L100, old=146499.72
$ws.Cells.Item(100, 14).Value = -129644.125  # N100, old=-147581.72
$ws.Cells.Item(126, 8).Value = 6197.4287  # H126, old=6229.2
$ws.Cells.Item(126, 9).Value = 4778.727  # I126, old=5386.625
$ws.Cells.Item(126, 10).Value = 11399.333  # J126, old=9599.5
$ws.Cells.Item(126, 11).Value = 14336.181  # K126, old=16159.875
$ws.Cells.Item(126, 12).Value = 34197.999  # L126, old=28798.5
$ws.Cells.Item(126, 13).Value = -11866.181  # M126, old=-13689.875
$ws.Cells.Item(126, 14).Value = -39137.999  # N126, old=-33738.5
$ws.Cells.Item(132, 8).Value = 8635.909  # H132, old=9091.846
$ws.Cells.Item(132, 9).Value = 8744.808999999999  # I132, old=9291.067999999999
$ws.Cells.Item(132, 11).Value = 26234.427  # K132, old=27873.204
$ws.Cells.Item(132, 13).Value = -23704.427  # M132, old=-25343.204
$ws.Cells.Item(136, 8).Value = 18185.174  # H136, old=19070.867
$ws.Cells.Item(136, 9).Value = 20760.904  # I136, old=21609.04
$ws.Cells.Item(136, 10).Value = 6009  # J136, old=6380
$ws.Cells.Item(136, 11).Value = 62282.712  # K136, old=64827.12
$ws.Cells.Item(136, 12).Value = 18027  # L136, old=19140
$ws.Cells.Item(136, 13).Value = -59732.712  # M136, old=-62277.12
$ws.Cells.Item(136, 14).Value = -23127  # N136, old=-24240

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(40, 8).Value = 33400  # H40, old=27625
$ws.Cells.Item(40, 9).Value = 30666.666  # I40, old=25166.666
$ws.Cells.Item(40, 10).Value = 37500  # J40, old=35000
$ws.Cells.Item(40, 11).Value = 30666.666  # K40, old=25166.666
$ws.Cells.Item(40, 12).Value = 37500  # L40, old=35000
$ws.Cells.Item(40, 13).Value = -30517.666  # M40, old=-25017.666
$ws.Cells.Item(40, 14).Value = -37798  # N40, old=-35298
$ws.Cells.Item(96, 8).Value = 2090  # H96, old=2112.75
$ws.Cells.Item(96, 10).Value = 2275  # J96, old=2367
$ws.Cells.Item(96, 12).Value = 2275  # L96, old=2367
$ws.Cells.Item(96, 14).Value = -5021  # N96, old=-5113
$ws.Cells.Item(107, 8).Value = 35720250  # H107, old=35720290
$ws.Cells.Item(107, 9).Value = 45457840  # I107, old=45457890
$ws.Cells.Item(107, 11).Value = 136373520  # K107, old=136373670
$ws.Cells.Item(107, 13).Value = -136371600  # M107, old=-136371750
$ws.Cells.Item(136, 8).Value = 5213.759  # H136, old=5747.346
$ws.Cells.Item(136, 9).Value = 5729.2173  # I136, old=6500.2
$ws.Cells.Item(136, 11).Value = 17187.6519  # K136, old=19500.6
$ws.Cells.Item(136, 13).Value = -14637.6519  # M136, old=-16950.6
